$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("G2").Value = 1.53
$ws.Range("H2").Value = 3.9
$ws.Range("I2").Value = 6.5
$ws.Range("J2").Value = 2.1
$ws.Range("K2").Value = 2.2
$ws.Range("AC2").Value = 8.5
$ws.Range("AD2").Value = 8
$ws.Range("AI2").Value = 34
$ws.Range("AN2").Value = 3.25
$ws.Range("AU2").Value = 9.5
$ws.Range("AV2").Value = 67
$ws.Range("AX2").Value = 41

# --- Row 3 ---
$ws.Range("O3").Value = 1.3
$ws.Range("P3").Value = 3.4
$ws.Range("Q3").Value = 2.03
$ws.Range("R3").Value = 1.83

# --- Row 5 ---
$ws.Range("H5").Value = 3
$ws.Range("I5").Value = 2.27
$ws.Range("K5").Value = 2.02
$ws.Range("L5").Value = 2.85
$ws.Range("N5").Value = 6.4
$ws.Range("O5").Value = 1.37
$ws.Range("P5").Value = 2.87
$ws.Range("T5").Value = 2.65
$ws.Range("U5").Value = 1.78
$ws.Range("V5").Value = 1.93
$ws.Range("W5").Value = 9.25
$ws.Range("X5").Value = 17.5
$ws.Range("AC5").Value = 6.4
$ws.Range("AD5").Value = 5.8
$ws.Range("AE5").Value = 13.5
$ws.Range("AF5").Value = 65
$ws.Range("AH5").Value = 7
$ws.Range("AI5").Value = 10.75
$ws.Range("AJ5").Value = 9
$ws.Range("AK5").Value = 23
$ws.Range("AL5").Value = 19.5
$ws.Range("AP5").Value = 23
$ws.Range("AT5").Value = 2.65
$ws.Range("AU5").Value = 6.6
$ws.Range("AW5").Value = 4.2
$ws.Range("AX5").Value = 12
$ws.Range("AY5").Value = 19.5
$ws.Range("AZ5").Value = 50
$ws.Range("BA5").Value = 80

# --- Row 8 ---
$ws.Range("H8").Value = 8
$ws.Range("I8").Value = 13
$ws.Range("J8").Value = 1.44
$ws.Range("L8").Value = 10
$ws.Range("N8").Value = 15
$ws.Range("Q8").Value = 1.33
$ws.Range("R8").Value = 3.25
$ws.Range("U8").Value = 2.1
$ws.Range("V8").Value = 1.67
$ws.Range("X8").Value = 7
$ws.Range("Y8").Value = 12
$ws.Range("Z8").Value = 7
$ws.Range("AD8").Value = 17
$ws.Range("AL8").Value = 81
$ws.Range("AM8").Value = 67
$ws.Range("AU8").Value = 11

# --- Row 10 ---
$ws.Range("G10").Value = 3.5
$ws.Range("H10").Value = 3
$ws.Range("K10").Value = 1.93
$ws.Range("L10").Value = 2.8
$ws.Range("N10").Value = 7.6
$ws.Range("O10").Value = 1.44
$ws.Range("P10").Value = 2.42
$ws.Range("Q10").Value = 2.25
$ws.Range("R10").Value = 1.5
$ws.Range("S10").Value = 1.5
$ws.Range("T10").Value = 2.25
$ws.Range("U10").Value = 1.98
$ws.Range("V10").Value = 1.65
$ws.Range("W10").Value = 8.25
$ws.Range("X10").Value = 17
$ws.Range("AB10").Value = 50
$ws.Range("AC10").Value = 7
$ws.Range("AD10").Value = 6
$ws.Range("AE10").Value = 17
$ws.Range("AF10").Value = 100
$ws.Range("AH10").Value = 6
$ws.Range("AI10").Value = 9
$ws.Range("AJ10").Value = 9
$ws.Range("AK10").Value = 20
$ws.Range("AL10").Value = 20
$ws.Range("AM10").Value = 37
$ws.Range("AO10").Value = 21
$ws.Range("AT10").Value = 2.22
$ws.Range("AU10").Value = 7.7
$ws.Range("AV10").Value = 90
$ws.Range("AY10").Value = 23
$ws.Range("BA10").Value = 100
$ws.Range("BB10").Value = 350
